$wb = $excel.ActiveWorkbook

# Rename the first sheet from "Sheet1" to "BOM" (auto-updates the
# localSheetId="0" defined names that reference Sheet1!...)
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "BOM"

# The BOM row for C12/C13 (10uF 1206 caps) was mislabeled; correct the
# value to "1uF".
$ws.Range("C8").Value = "1uF"

# Move the active cell/selection back up to D7 (was left at D29).
$ws.Range("D7").Select()
